$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Income"

# Header row + 13 income entries (Icon, Source, Amount, Date)
$incomeRows = @(
    @("Icon", "Source", "Amount", "Date"),
    @("📊", "Consulting", 20000, "24 Jul 2025"),
    @("📝", "Content Writing", 6200, "23 Jul 2025"),
    @("📦", "E-commerce Sales", 9500, "22 Jul 2025"),
    @("💰", "Side Hustle", 5200, "21 Jul 2025"),
    @("💻", "Web Development", 15000, "20 Jul 2025"),
    @("🎮", "Game Development", 13000, "19 Jul 2025"),
    @("😃", "freelance", 15000, "19 Jul 2025"),
    @("💰", "Investments", 9500, "18 Jul 2025"),
    @("🎓", "Online Tutoring", 7400, "18 Jul 2025"),
    @("📸", "Photography", 11200, "17 Jul 2025"),
    @("💰", "Bonus", 8000, "15 Jul 2025"),
    @("🏟️", "salary", 40000, "15 Jul 2025"),
    @("💰", "Salary", 45000, "01 Jul 2025")
)

for ($i = 0; $i -lt $incomeRows.Count; $i++) {
    $r = $i + 1
    $row = $incomeRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
